$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (Emily Jacobson): add Favorite Ice Cream / Favorite Pizza Toppings
$ws.Range("E6").Value = "mint chocolate chip"
$ws.Range("F6").Value = "onion"

# Row 12 (Lukas Larson): add Favorite Ice Cream / Favorite Pizza Toppings
$ws.Range("E12").Value = "I'm Lactose intolerant"
$ws.Range("F12").Value = "Still Lactose Intolerant"

# Row 13 (Luke Bertram): clear previous Favorite Ice Cream / Favorite Pizza Toppings
$ws.Range("E13").ClearContents()
$ws.Range("F13").ClearContents()
